$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3741247334384212
$ws.Range("D2").Value = 0.04649824640627287
$ws.Range("E2").Value = 0.1456202660292512
$ws.Range("F2").Value = 1.108027689185711
$ws.Range("G2").Value = 0.00246990454627088
$ws.Range("K2").Value = 1.13387091036364
$ws.Range("L2").Value = 0.1187981008196708
$ws.Range("M2").Value = 0.414265231808109
$ws.Range("N2").Value = 1.432106154888601
$ws.Range("O2").Value = 3.976110968011028
$ws.Range("C3").Value = 0.3704893607892217
$ws.Range("D3").Value = 0.0456094257289763
$ws.Range("E3").Value = 0.1453094501481331
$ws.Range("F3").Value = 1.105147322339121
$ws.Range("G3").Value = 0.002473073296272595
$ws.Range("K3").Value = 1.025899616534019
$ws.Range("L3").Value = 0.1192635899362919
$ws.Range("M3").Value = 0.3906776901091931
$ws.Range("N3").Value = 1.44470177941929
$ws.Range("O3").Value = 3.982991701570995
$ws.Range("C4").Value = 0.3684495232528775
$ws.Range("D4").Value = 0.04505973823963672
$ws.Range("E4").Value = 0.1451902544362795
$ws.Range("F4").Value = 1.104085456885372
$ws.Range("G4").Value = 0.002475123370083089
$ws.Range("K4").Value = 0.9597441310163219
$ws.Range("L4").Value = 0.119603175217744
$ws.Range("M4").Value = 0.3763397990846826
$ws.Range("N4").Value = 1.45297337215154
$ws.Range("O4").Value = 3.98965418983525
$ws.Range("C5").Value = 0.3676666684583552
$ws.Range("D5").Value = 0.04483476028687505
$ws.Range("E5").Value = 0.1451597047179476
$ws.Range("F5").Value = 1.103830273740535
$ws.Range("G5").Value = 0.002475985136549722
$ws.Range("K5").Value = 0.9328216157491624
$ws.Range("L5").Value = 0.1197550820984468
$ws.Range("M5").Value = 0.3705337171827168
$ws.Range("N5").Value = 1.456479448064947
$ws.Range("O5").Value = 3.992981566630618
$ws.Range("C6").Value = 0.3675396002182367
$ws.Range("D6").Value = 0.04479734443001604
$ws.Range("E6").Value = 0.1451557207588223
$ws.Range("F6").Value = 1.103798618615876
$ws.Range("G6").Value = 0.002476129825681804
$ws.Range("K6").Value = 0.928353389823144
$ws.Range("L6").Value = 0.1197811229298189
$ws.Range("M6").Value = 0.3695718479469008
$ws.Range("N6").Value = 1.457069806031292
$ws.Range("O6").Value = 3.993571041986257
$ws.Range("C7").Value = 0.3684387693907922
$ws.Range("D7").Value = 0.04505670803699502
$ws.Range("E7").Value = 0.1451897694454622
$ws.Range("F7").Value = 1.104081296785708
$ws.Range("G7").Value = 0.002475134885486963
$ws.Range("K7").Value = 0.9593808954599865
$ws.Range("L7").Value = 0.1196051691306366
$ws.Range("M7").Value = 0.3762613471696739
$ws.Range("N7").Value = 1.453020108227996
$ws.Range("O7").Value = 3.989696585414634
$ws.Range("C8").Value = 0.372831368417252
$ws.Range("D8").Value = 0.0461926103013468
$ws.Range("E8").Value = 0.145498230418255
$ws.Range("F8").Value = 1.106887759594414
$ws.Range("G8").Value = 0.002470975501157303
$ws.Range("K8").Value = 1.096614264246483
$ws.Range("L8").Value = 0.1189474415523577
$ws.Range("M8").Value = 0.4061023186759272
$ws.Range("N8").Value = 1.436337555164108
$ws.Range("O8").Value = 3.977977119856661
$ws.Range("C9").Value = 0.3829701495823485
$ws.Range("D9").Value = 0.04838814519739287
$ws.Range("E9").Value = 0.1466714996726672
$ws.Range("F9").Value = 1.118008068559362
$ws.Range("G9").Value = 0.002463644055040848
$ws.Range("K9").Value = 1.366788506102012
$ws.Range("L9").Value = 0.118084322554818
$ws.Range("M9").Value = 0.4657617337626405
$ws.Range("N9").Value = 1.407886479943528
$ws.Range("O9").Value = 3.974369993003876
$ws.Range("C10").Value = 0.3913489466884812
$ws.Range("D10").Value = 0.04998100805393335
$ws.Range("E10").Value = 0.147880186480041
$ws.Range("F10").Value = 1.129618467228056
$ws.Range("G10").Value = 0.002458755514039978
$ws.Range("K10").Value = 1.565890928307567
$ws.Range("L10").Value = 0.1177104503688469
$ws.Range("M10").Value = 0.5102823541263035
$ws.Range("N10").Value = 1.389577120850305
$ws.Range("O10").Value = 3.983583488554444
$ws.Range("C11").Value = 0.395362724855147
$ws.Range("D11").Value = 0.0507011175196439
$ws.Range("E11").Value = 0.1485053937464933
$ws.Range("F11").Value = 1.135651035377933
$ws.Range("G11").Value = 0.002456638607492169
$ws.Range("K11").Value = 1.656591896664793
$ws.Range("L11").Value = 0.1175969216881541
$ws.Range("M11").Value = 0.5306843573468711
$ws.Range("N11").Value = 1.381809876857893
$ws.Range("O11").Value = 3.990362812851089
$ws.Range("C12").Value = 0.3969117009916943
$ws.Range("D12").Value = 0.05097314345271542
$ws.Range("E12").Value = 0.1487529809054244
$ws.Range("F12").Value = 1.138043639226609
$ws.Range("G12").Value = 0.002455852281707313
$ws.Range("K12").Value = 1.690955370944494
$ws.Range("L12").Value = 0.1175620648552389
$ws.Range("M12").Value = 0.538431336965786
$ws.Range("N12").Value = 1.378949347694487
$ws.Range("O12").Value = 3.993302968172372
$ws.Range("C13").Value = 0.3965768099955653
$ws.Range("D13").Value = 0.05091458761226164
$ws.Range("E13").Value = 0.1486991767545902
$ws.Range("F13").Value = 1.137523533441126
$ws.Range("G13").Value = 0.002456020951647625
$ws.Range("K13").Value = 1.683553849643317
$ws.Range("L13").Value = 0.1175692100935422
$ws.Range("M13").Value = 0.5367619487839761
$ws.Range("N13").Value = 1.379561822452047
$ws.Range("O13").Value = 3.992653152100644
$ws.Range("C14").Value = 0.3954895781326968
$ws.Range("D14").Value = 0.05072351066515068
$ws.Range("E14").Value = 0.1485255457653984
$ws.Range("F14").Value = 1.135845706579175
$ws.Range("G14").Value = 0.002456573609719324
$ws.Range("K14").Value = 1.659418668292176
$ws.Range("L14").Value = 0.1175938909628869
$ws.Range("M14").Value = 0.5313212831393201
$ws.Range("N14").Value = 1.381572921022787
$ws.Range("O14").Value = 3.990597221613996
$ws.Range("C15").Value = 0.3948273988026756
$ws.Range("D15").Value = 0.05060638352758673
$ws.Range("E15").Value = 0.1484206026770671
$ws.Range("F15").Value = 1.134832086617251
$ws.Range("G15").Value = 0.002456914119195958
$ws.Range("K15").Value = 1.644637341674468
$ws.Range("L15").Value = 0.1176100680739154
$ws.Range("M15").Value = 0.5279914684285529
$ws.Range("N15").Value = 1.382815293342688
$ws.Range("O15").Value = 3.989386500998165
$ws.Range("C16").Value = 0.3910907031673219
$ws.Range("D16").Value = 0.04993385527959759
$ws.Range("E16").Value = 0.1478408440569545
$ws.Range("F16").Value = 1.129239354081037
$ws.Range("G16").Value = 0.002458896003904015
$ws.Range("K16").Value = 1.559965863037576
$ws.Range("L16").Value = 0.1177190078318837
$ws.Range("M16").Value = 0.5089520154989131
$ws.Range("N16").Value = 1.390096030591742
$ws.Range("O16").Value = 3.983192591532031
$ws.Range("C17").Value = 0.3888501310774188
$ws.Range("D17").Value = 0.04952011768523334
$ws.Range("E17").Value = 0.1475044824787588
$ws.Range("F17").Value = 1.126000881787988
$ws.Range("G17").Value = 0.002460139157081976
$ws.Range("K17").Value = 1.508054467597049
$ws.Range("L17").Value = 0.1178003238171925
$ws.Range("M17").Value = 0.4973099713339906
$ws.Range("N17").Value = 1.394706391511342
$ws.Range("O17").Value = 3.980056270028541
$ws.Range("C18").Value = 0.3875804493032149
$ws.Range("D18").Value = 0.04928172503573336
$ws.Range("E18").Value = 0.1473181105279942
$ws.Range("F18").Value = 1.12420886716518
$ws.Range("G18").Value = 0.002460864253898991
$ws.Range("K18").Value = 1.478208566965634
$ws.Range("L18").Value = 0.1178524169948183
$ws.Range("M18").Value = 0.4906278503628201
$ws.Range("N18").Value = 1.397411027228706
$ws.Range("O18").Value = 3.978495892418749
$ws.Range("C19").Value = 0.3871538273602937
$ws.Range("D19").Value = 0.04920093749821319
$ws.Range("E19").Value = 0.1472562267491462
$ws.Range("F19").Value = 1.123614252667537
$ws.Range("G19").Value = 0.002461111490530347
$ws.Range("K19").Value = 1.468105391093445
$ws.Range("L19").Value = 0.1178709688975701
$ws.Range("M19").Value = 0.488367825329945
$ws.Range("N19").Value = 1.398335853132181
$ws.Range("O19").Value = 3.978009380889119
$ws.Range("C20").Value = 0.3890866738888406
$ws.Range("D20").Value = 0.04956420453593324
$ws.Range("E20").Value = 0.1475395545309013
$ws.Range("F20").Value = 1.126338307143428
$ws.Range("G20").Value = 0.002460005779981872
$ws.Range("K20").Value = 1.513579276237749
$ws.Range("L20").Value = 0.1177911167279611
$ws.Range("M20").Value = 0.4985478322325534
$ws.Range("N20").Value = 1.394210138236993
$ws.Range("O20").Value = 3.98036492424373
$ws.Range("C21").Value = 0.3958081364681334
$ws.Range("D21").Value = 0.05077965275640395
$ws.Range("E21").Value = 0.1485762513700273
$ws.Range("F21").Value = 1.136335586718644
$ws.Range("G21").Value = 0.002456410865927374
$ws.Range("K21").Value = 1.666507306711139
$ws.Range("L21").Value = 0.1175864208179398
$ws.Range("M21").Value = 0.5329187659263397
$ws.Range("N21").Value = 1.380980021055564
$ws.Range("O21").Value = 3.991190969381393
$ws.Range("C22").Value = 0.4003702844961765
$ws.Range("D22").Value = 0.05157014050841724
$ws.Range("E22").Value = 0.1493169443654949
$ws.Range("F22").Value = 1.143500167248348
$ws.Range("G22").Value = 0.00245415053304823
$ws.Range("K22").Value = 1.766552949330674
$ws.Range("L22").Value = 0.1175000527324066
$ws.Range("M22").Value = 0.5555054984154992
$ws.Range("N22").Value = 1.372804084116829
$ws.Range("O22").Value = 4.000440763051159
$ws.Range("C23").Value = 0.3979199018367865
$ws.Range("D23").Value = 0.05114860305248214
$ws.Range("E23").Value = 0.1489158452543045
$ws.Range("F23").Value = 1.139618507556463
$ws.Range("G23").Value = 0.002455348782038345
$ws.Range("K23").Value = 1.713148190313916
$ws.Range("L23").Value = 0.1175418098871326
$ws.Range("M23").Value = 0.543439344884078
$ws.Range("N23").Value = 1.377124672337729
$ws.Range("O23").Value = 3.995304751960248
$ws.Range("C24").Value = 0.388979675403732
$ws.Range("D24").Value = 0.04954427451731647
$ws.Range("E24").Value = 0.1475236766306871
$ws.Range("F24").Value = 1.126185539648688
$ws.Range("G24").Value = 0.002460066047487859
$ws.Range("K24").Value = 1.511081514340503
$ws.Range("L24").Value = 0.1177952626087269
$ws.Range("M24").Value = 0.4979881608938754
$ws.Range("N24").Value = 1.394434325885427
$ws.Range("O24").Value = 3.980224625611186
$ws.Range("C25").Value = 0.3800640682825076
$ws.Range("D25").Value = 0.04779769086022156
$ws.Range("E25").Value = 0.146293231610084
$ws.Range("F25").Value = 1.11439674317073
$ws.Range("G25").Value = 0.002465539605488787
$ws.Range("K25").Value = 1.293589975104908
$ws.Range("L25").Value = 0.1182721253845465
$ws.Range("M25").Value = 0.4495008147406594
$ws.Range("N25").Value = 1.415127410429584
$ws.Range("O25").Value = 3.973266226057234
